$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet (per-fund holdings), placed right after
#    "2021-Q4" and before "总计". We duplicate "2021-Q4" so the new sheet
#    inherits the exact same layout / cell styles (header row style, index
#    column style, sheetPr, pageMargins, etc.), then we overwrite its data.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)
$newQSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newQSheet.Name = "2022-Q1"

# "2021-Q4" has 3 data rows (rows 2-4); the new sheet only needs 2 (rows 2-3),
# so drop the extra row.
$newQSheet.Rows.Item(4).Delete()

# Row 2: 九泰天宝灵活配置混合A
$newQSheet.Cells.Item(2,1).Value = 0
$newQSheet.Cells.Item(2,2).Value = "'000892"
$newQSheet.Cells.Item(2,2).Style = "Normal"
$newQSheet.Cells.Item(2,3).Value = "九泰天宝灵活配置混合A"
$newQSheet.Cells.Item(2,4).Value = "'0.07"
$newQSheet.Cells.Item(2,4).Style = "Normal"
$newQSheet.Cells.Item(2,5).Value = "'90.81"
$newQSheet.Cells.Item(2,5).Style = "Normal"
$newQSheet.Cells.Item(2,6).Value = "'4.56"
$newQSheet.Cells.Item(2,6).Style = "Normal"
$newQSheet.Cells.Item(2,7).Value = "'0.0032"
$newQSheet.Cells.Item(2,7).Style = "Normal"
$newQSheet.Cells.Item(2,8).Value = 8

# Row 3: 九泰天宝灵活配置混合C
$newQSheet.Cells.Item(3,1).Value = 1
$newQSheet.Cells.Item(3,2).Value = "'002028"
$newQSheet.Cells.Item(3,2).Style = "Normal"
$newQSheet.Cells.Item(3,3).Value = "九泰天宝灵活配置混合C"
$newQSheet.Cells.Item(3,4).Value = "'0.00"
$newQSheet.Cells.Item(3,4).Style = "Normal"
$newQSheet.Cells.Item(3,5).Value = "'90.81"
$newQSheet.Cells.Item(3,5).Style = "Normal"
$newQSheet.Cells.Item(3,6).Value = "'4.56"
$newQSheet.Cells.Item(3,6).Style = "Normal"
$newQSheet.Cells.Item(3,7).Value = 0
$newQSheet.Cells.Item(3,8).Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new top data row for the
#    2022-Q1 quarter (2 funds held, 0 亿元 market value), pushing the
#    existing quarters down, and renumber the index column (A) 0..4.
# ---------------------------------------------------------------------------
$sumSheet = $wb.Worksheets.Item("总计")
$sumSheet.Rows.Item(2).Insert()

# The blank row Excel creates on insert picks up borders/bold from the
# header row above it - reset B2:D2 back to the plain "Normal" style used
# by the other data rows.
$sumSheet.Range("B2:D2").Style = "Normal"

# Copy the index-column style (bold + border) from A3 onto the new A2 cell.
$sumSheet.Range("A3").Copy()
$sumSheet.Range("A2").PasteSpecial(-4122)

$sumSheet.Cells.Item(2,1).Value = 0
$sumSheet.Cells.Item(2,2).Value = "2022-Q1"
$sumSheet.Cells.Item(2,3).Value = 2
$sumSheet.Cells.Item(2,4).Value = 0

$sumSheet.Cells.Item(3,1).Value = 1
$sumSheet.Cells.Item(4,1).Value = 2
$sumSheet.Cells.Item(5,1).Value = 3
$sumSheet.Cells.Item(6,1).Value = 4
